$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.812.78'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').Value = '1.766.52'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.47%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9991'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4258'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3607'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.24'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07440'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.102'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9936'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.109'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.290'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.80%  '
$ws.Range('D16').Value = '1.792.66'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.05'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001060'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06383'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9980'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.956'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.40%  '
$ws.Range('D23').Value = '27.841.50'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.159'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.55%  '
$ws.Range('D28').Value = '1.992.27'
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.153'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.87'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.168'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.687'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09000'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.67'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02310'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2113'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.21%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.053'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06078'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6401'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.179'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.864'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.391'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.66'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5956'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.694'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.984'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.144'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06888'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.11%  '
